$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Sheet2 - Numbers"): add a new column AA with values 100-129 ---
$ws2 = $wb.Worksheets.Item("Sheet2 - Numbers")

for ($row = 1; $row -le 30; $row++) {
    $ws2.Cells.Item($row, 27).Value = 99 + $row
}

# Make Sheet2 the active sheet/tab and select AA1:AA30 (new data column)
$ws2.Activate()
$ws2.Range("AA1:AA30").Select()

# --- Sheet4 ("Sheet4 - Dates"): switch page size to A4 ---
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9
